# Converts an "RRGGBB" hex color string into the BGR-packed integer
# expected by the PowerPoint COM `RGB` color property
# (VBA convention: value = R + G*256 + B*65536).
function HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 6: switch to the built-in "Medium Style 2 -
#        Accent 1" table style. ---------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{532D129F-856A-48F6-82E0-E73179350A86}")
    }
}

# --- 2. Switch the presentation's colour theme from "Integral" back
#        to the default Office colours. ---------------------------------
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $themeColors.Item($i + 1).RGB = HexToRgbInt($officeColors[$i])
}
